$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# optimization_parameters sheet (7th tab) restructuring
# ---------------------------------------------------------------------------
$wsOpt = $wb.Worksheets.Item(7)

# Remove the obsolete "Deletion" row (old row 16) first so row numbers below
# it close up before we insert the new row above.
$wsOpt.Rows.Item(16).Delete()

# Insert a new blank row right after the "Model"/"production_function" row
# (old row 8) to host the new "L_curve" parameter.
$wsOpt.Rows.Item(9).Insert()

# The header row had "value" duplicated across C1:F1 -- drop those, leaving
# just the A1/B1 header pair.
$wsOpt.Range("C1:F1").ClearContents()

# Rename the "Model" parameter label to "production_function".
$wsOpt.Range("A8").Value = "production_function"

# Populate the newly inserted row with the "L_curve" parameter.
$wsOpt.Range("A9").Value = "L_curve"
$wsOpt.Range("B9").Value = 0
$wsOpt.Range("B9").NumberFormat = "0.00E+00"

# ---------------------------------------------------------------------------
# View state: the active tab & selection move from network_weights (tab 6)
# to optimization_parameters (tab 7), and the selection there becomes the
# cleared header cells C1:F1.
# ---------------------------------------------------------------------------
$wsOpt.Activate()
$wsOpt.Range("C1:F1").Select()
